$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2296.5
$ws.Range("I40").Value = 2052.4285
$ws.Range("J40").Value = 2638.2
$ws.Range("K40").Value = 2052.4285
$ws.Range("L40").Value = 2638.2
$ws.Range("M40").Value = -1877.4285
$ws.Range("N40").Value = -2988.2
$ws.Range("H75").Value = 28545
$ws.Range("J75").Value = 24597
$ws.Range("L75").Value = 24597
$ws.Range("N75").Value = -26469
$ws.Range("H78").Value = 28545
$ws.Range("J78").Value = 24597
$ws.Range("L78").Value = 73791
$ws.Range("N78").Value = -83151
$ws.Range("H96").Value = 1127.2667
$ws.Range("I96").Value = 1625.8889
$ws.Range("J96").Value = 379.33334
$ws.Range("K96").Value = 4877.6667
$ws.Range("L96").Value = 1138.00002
$ws.Range("M96").Value = -3504.6667
$ws.Range("N96").Value = -3884.00002
$ws.Range("H103").Value = 71429704
$ws.Range("I103").Value = 100000790
$ws.Range("K103").Value = 300002370
$ws.Range("M103").Value = -300001784
$ws.Range("H137").Value = 2990
$ws.Range("I137").Value = 2990
$ws.Range("K137").Value = 8970
$ws.Range("M137").Value = -6420

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1576
$ws.Range("I2").Value = 1514.4445
$ws.Range("J2").Value = 1797.6
$ws.Range("K2").Value = 1514.4445
$ws.Range("L2").Value = 1797.6
$ws.Range("M2").Value = -1401.4445
$ws.Range("N2").Value = -2023.6
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("H32").Value = 1013.8605
$ws.Range("I32").Value = 1037
$ws.Range("K32").Value = 1037
$ws.Range("M32").Value = -750
$ws.Range("H45").Value = 614
$ws.Range("J45").Value = 614
$ws.Range("L45").Value = 614
$ws.Range("N45").Value = -1368
$ws.Range("H63").Value = 3543.139
$ws.Range("I63").Value = 1951.9667
$ws.Range("K63").Value = 1951.9667
$ws.Range("M63").Value = -1265.9667
$ws.Range("H66").Value = 3543.139
$ws.Range("I66").Value = 1951.9667
$ws.Range("K66").Value = 9759.833500000001
$ws.Range("M66").Value = -6327.833500000001
$ws.Range("H74").Value = 1856.826
$ws.Range("I74").Value = 1564.2106
$ws.Range("J74").Value = 3246.75
$ws.Range("K74").Value = 1564.2106
$ws.Range("L74").Value = 3246.75
$ws.Range("M74").Value = -690.2106000000001
$ws.Range("N74").Value = -4994.75
$ws.Range("H77").Value = 1856.826
$ws.Range("I77").Value = 1564.2106
$ws.Range("J77").Value = 3246.75
$ws.Range("K77").Value = 7821.053000000001
$ws.Range("L77").Value = 16233.75
$ws.Range("M77").Value = -3453.053000000001
$ws.Range("N77").Value = -24969.75
$ws.Range("H97").Value = 521.8
$ws.Range("I97").Value = 521.8
$ws.Range("K97").Value = 521.8
$ws.Range("M97").Value = -25.79999999999995
$ws.Range("H116").Value = 1576
$ws.Range("I116").Value = 1514.4445
$ws.Range("J116").Value = 1797.6
$ws.Range("K116").Value = 1514.4445
$ws.Range("L116").Value = 1797.6
$ws.Range("M116").Value = 779.5554999999999
$ws.Range("N116").Value = -6385.6
$ws.Range("H122").Value = 865.3
$ws.Range("I122").Value = 865.3
$ws.Range("K122").Value = 2595.9
$ws.Range("M122").Value = -145.8999999999996

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1576
$ws.Range("I3").Value = 1514.4445
$ws.Range("J3").Value = 1797.6
$ws.Range("K3").Value = 1514.4445
$ws.Range("L3").Value = 1797.6
$ws.Range("M3").Value = -1400.4445
$ws.Range("N3").Value = -2025.6
$ws.Range("H82").Value = 7280.875
$ws.Range("I82").Value = 7280.875
$ws.Range("K82").Value = 7280.875
$ws.Range("M82").Value = -6897.875
$ws.Range("H85").Value = 7280.875
$ws.Range("I85").Value = 7280.875
$ws.Range("K85").Value = 7280.875
$ws.Range("M85").Value = -5954.875
$ws.Range("H94").Value = 381.2857
$ws.Range("I94").Value = 381.2857
$ws.Range("K94").Value = 381.2857
$ws.Range("M94").Value = 69.71429999999998
$ws.Range("H105").Value = 2839
$ws.Range("I105").Value = 1756.6364
$ws.Range("K105").Value = 1756.6364
$ws.Range("M105").Value = -9.636400000000094
$ws.Range("H107").Value = 1892.5333
$ws.Range("I107").Value = 1906.2858
$ws.Range("K107").Value = 1906.2858
$ws.Range("M107").Value = 13.71419999999989

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1563.5333
$ws.Range("I22").Value = 1077.6
$ws.Range("J22").Value = 1806.5
$ws.Range("K22").Value = 1077.6
$ws.Range("L22").Value = 1806.5
$ws.Range("M22").Value = -727.5999999999999
$ws.Range("N22").Value = -2506.5
$ws.Range("H31").Value = 4524.4
$ws.Range("I31").Value = 3655.625
$ws.Range("K31").Value = 3655.625
$ws.Range("M31").Value = -3360.625
$ws.Range("H34").Value = 4524.4
$ws.Range("I34").Value = 3655.625
$ws.Range("K34").Value = 3655.625
$ws.Range("M34").Value = -3453.625
$ws.Range("H51").Value = 55048.5
$ws.Range("J51").Value = 55048.5
$ws.Range("L51").Value = 55048.5
$ws.Range("N51").Value = -56520.5
$ws.Range("H58").Value = 3014.5334
$ws.Range("I58").Value = 2507.2727
$ws.Range("K58").Value = 2507.2727
$ws.Range("M58").Value = -2304.2727
$ws.Range("H61").Value = 55048.5
$ws.Range("J61").Value = 55048.5
$ws.Range("L61").Value = 55048.5
$ws.Range("N61").Value = -55744.5
$ws.Range("H107").Value = 1629.0385
$ws.Range("J107").Value = 1842
$ws.Range("L107").Value = 1842
$ws.Range("N107").Value = -5682
$ws.Range("H122").Value = 1300
$ws.Range("I122").Value = 1300
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3900
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -1450
$ws.Range("H136").Value = 3014.5334
$ws.Range("I136").Value = 2507.2727
$ws.Range("K136").Value = 7521.8181
$ws.Range("M136").Value = -4971.8181

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 47.666668
$ws.Range("I6").Value = 47.666668
$ws.Range("K6").Value = 143.000004
$ws.Range("M6").Value = -30.00000399999999
$ws.Range("H14").Value = 1160.7858
$ws.Range("I14").Value = 1160.7858
$ws.Range("K14").Value = 3482.3574
$ws.Range("M14").Value = -3309.3574
$ws.Range("H68").Value = 6000
$ws.Range("J68").Value = 6000
$ws.Range("L68").Value = 18000
$ws.Range("N68").Value = -19622
$ws.Range("H71").Value = 6000
$ws.Range("J71").Value = 6000
$ws.Range("L71").Value = 54000
$ws.Range("N71").Value = -62112
$ws.Range("H80").Value = 2874.5
$ws.Range("J80").Value = 2499.5
$ws.Range("L80").Value = 7498.5
$ws.Range("N80").Value = -9370.5
$ws.Range("H83").Value = 2874.5
$ws.Range("J83").Value = 2499.5
$ws.Range("L83").Value = 22495.5
$ws.Range("N83").Value = -31855.5
$ws.Range("H97").Value = 1545.909
$ws.Range("I97").Value = 360.8
$ws.Range("J97").Value = 2533.5
$ws.Range("K97").Value = 1082.4
$ws.Range("L97").Value = 7600.5
$ws.Range("M97").Value = -586.4000000000001
$ws.Range("N97").Value = -8592.5
$ws.Range("H136").Value = 2359.8
$ws.Range("I136").Value = 2359.8
$ws.Range("K136").Value = 7079.400000000001
$ws.Range("M136").Value = -1979.400000000001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 6332.6665
$ws.Range("I5").Value = 6332.6665
$ws.Range("K5").Value = 6332.6665
$ws.Range("M5").Value = -6220.6665
$ws.Range("H102").Value = 2038.6
$ws.Range("I102").Value = 2042.8889
$ws.Range("K102").Value = 2042.8889
$ws.Range("M102").Value = -420.8888999999999
$ws.Range("H107").Value = 2116.5881
$ws.Range("I107").Value = 894.1818
$ws.Range("K107").Value = 894.1818
$ws.Range("M107").Value = 1025.8182

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 7800
$ws.Range("I16").Value = 1700
$ws.Range("J16").Value = 20000
$ws.Range("K16").Value = 1700
$ws.Range("L16").Value = 20000
$ws.Range("M16").Value = -1530
$ws.Range("N16").Value = -20340
$ws.Range("H26").Value = 5000
$ws.Range("J26").Value = 5000
$ws.Range("L26").Value = 5000
$ws.Range("N26").Value = -5590
$ws.Range("H132").Value = 2601
$ws.Range("I132").Value = 2500.375
$ws.Range("K132").Value = 7501.125
$ws.Range("M132").Value = -4971.125

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2146.4375
$ws.Range("J81").Value = 2948.5
$ws.Range("L81").Value = 5897
$ws.Range("N81").Value = -8019
$ws.Range("H84").Value = 2146.4375
$ws.Range("J84").Value = 2948.5
$ws.Range("L84").Value = 29485
$ws.Range("N84").Value = -40093
$ws.Range("H126").Value = 2199.4443
$ws.Range("I126").Value = 2224.375
$ws.Range("K126").Value = 6673.125
$ws.Range("M126").Value = -4203.125
